$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column "want-to-go" counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 13965
$ws1.Range("F3").Value = 329
$ws1.Range("F4").Value = 674
$ws1.Range("F5").Value = 237
$ws1.Range("F6").Value = 517
$ws1.Range("F7").Value = 1447
$ws1.Range("F8").Value = 137

# Sheet "全部类型" (sheet4): update F column "want-to-go" counts for matching rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 13965
$ws4.Range("F3").Value = 329
$ws4.Range("F4").Value = 674
$ws4.Range("F5").Value = 237
$ws4.Range("F8").Value = 517
$ws4.Range("F9").Value = 1447
$ws4.Range("F11").Value = 137
